# Week 15 simulations added.
# Appends newly simulated per-game numbers to the long space-separated
# number lists on the YDS and ST sheets, and updates the season-total
# numeric cells on OFF / DEF / ST / TURNS / PEN that roll those lists up.

$wb = $excel.ActiveWorkbook

function Append-Numbers {
    param(
        [string]$SheetName,
        [string]$CellAddr,
        [string]$Suffix
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $cell = $ws.Range($CellAddr)
    $cell.Value2 = [string]$cell.Value2 + $Suffix
}

# --- YDS sheet: append this week's rush/pass yardage-log entries ---
Append-Numbers "YDS" "B2" " 2 6 3 17 0 5 4 4 6 8 2 5 3 4 3 2 -5 3 1 1 -1 2 4 5 6 1 1 12 0 2 -1"
Append-Numbers "YDS" "B3" " 9 6 11 15 1 16 9 12 5 4 8 11 4 10 24 9 5 13 7 12"
Append-Numbers "YDS" "C2" " -4 3 5 2 2 -1 0 1"
Append-Numbers "YDS" "C3" " 8 14 3 7 5 7 3 12 3 31 7 14 6 5 5 9 9 31 41 6 6 3 2"

# --- ST sheet: append this week's special-teams-log entries ---
Append-Numbers "ST" "B4" " 66"
Append-Numbers "ST" "B5" " 16"
Append-Numbers "ST" "B6" " 25"
Append-Numbers "ST" "D3" " 35 37 35 39 46"
Append-Numbers "ST" "D4" " 0 0 0 0 0"
Append-Numbers "ST" "D5" " 2 0 0 1 10"

# --- OFF sheet: season totals updated for the new week ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value2 = 9
$ws.Range("C2").Value2 = 459
$ws.Range("D2").Value2 = 28
$ws.Range("F2").Value2 = 142
$ws.Range("G2").Value2 = 107
$ws.Range("I2").Value2 = 6
$ws.Range("J2").Value2 = 58
$ws.Range("L2").Value2 = 499
$ws.Range("M2").Value2 = 331
$ws.Range("Q2").Value2 = 982
$ws.Range("C3").Value2 = 237
$ws.Range("E3").Value2 = 57
$ws.Range("F3").Value2 = 177
$ws.Range("G3").Value2 = 63
$ws.Range("H3").Value2 = 45
$ws.Range("I3").Value2 = 94
$ws.Range("J3").Value2 = 107
$ws.Range("N3").Value2 = 32

# --- DEF sheet: season totals updated for the new week ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value2 = 324
$ws.Range("F2").Value2 = 105
$ws.Range("G2").Value2 = 115
$ws.Range("I2").Value2 = 16
$ws.Range("L2").Value2 = 610
$ws.Range("M2").Value2 = 405
$ws.Range("O2").Value2 = 43
$ws.Range("P2").Value2 = 26
$ws.Range("Q2").Value2 = 1041
$ws.Range("B3").Value2 = 23
$ws.Range("C3").Value2 = 356
$ws.Range("E3").Value2 = 52
$ws.Range("F3").Value2 = 207
$ws.Range("G3").Value2 = 59
$ws.Range("H3").Value2 = 40
$ws.Range("I3").Value2 = 117
$ws.Range("J3").Value2 = 107
$ws.Range("N3").Value2 = 30

# --- ST sheet: season totals updated for the new week ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value2 = 165
$ws.Range("D2").Value2 = 95
$ws.Range("F2").Value2 = 259
$ws.Range("G2").Value2 = 246
$ws.Range("J2").Value2 = 120
$ws.Range("K2").Value2 = 113
$ws.Range("L2").Value2 = 80
$ws.Range("M2").Value2 = 62
$ws.Range("B3").Value2 = 94

# --- TURNS sheet: season totals updated for the new week ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value2 = 15
$ws.Range("D3").Value2 = 12

# --- PEN sheet: season totals updated for the new week ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B3").Value2 = 23
